$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("B3").Value = -0.579113769187642
$ws.Range("C3").Value = 1.63506778697617
$ws.Range("D3").Value = 13.14160213705591
$ws.Range("E3").Value = 3.625134775019531
$ws.Range("F3").Value = 3.662792626675848
$ws.Range("G3").Value = 22

# Row 4 updates
$ws.Range("B4").Value = 0.6174351661100629
$ws.Range("C4").Value = 1.151642498604435
$ws.Range("D4").Value = 6.041101927293427
$ws.Range("E4").Value = 2.45786531919335
$ws.Range("F4").Value = 2.437800141539759
$ws.Range("G4").Value = 21
